# Weekly update: a new price observation is added to the top of the
# historical list (row 9), pushing every existing observation down by
# one row. The oldest observation (previously on row 26) ends up on the
# new row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; Excel shifts rows 9:26 down to 10:27
# and the sheet dimension grows from A1:R26 to A1:R27 automatically.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new weekly observation.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44575
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112028
$ws.Range("G9").Value = "Sandia"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 380
$ws.Range("L9").Value = 400
$ws.Range("M9").Value = 390
$ws.Range("N9").Value = "$/kilo (volumen en unidades)"
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 390
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
